$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.640.03'
$ws.Range('E2').Value = '  +4.13%  '
$ws.Range('D3').Value = '2.257.84'
$ws.Range('E3').Value = '  +2.68%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '91.09'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.74%  '
$ws.Range('E7').Value = '  +3.43%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.479'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.05'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.14'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.31%  '
$ws.Range('E12').Value = '  +2.08%  '
$ws.Range('E13').Value = '  +1.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.56'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.03%  '
$ws.Range('D15').Value = '2.607.72'
$ws.Range('E15').Value = '  +2.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.11'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.79%  '
$ws.Range('D17').Value = '2.251.04'
$ws.Range('E17').Value = '  -0.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.757'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.36%  '
$ws.Range('D19').Value = '41.563.79'
$ws.Range('E19').Value = '  +4.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.45'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +11.55%  '
$ws.Range('E21').Value = '  +1.96%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.89'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.50'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '240.25'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.22%  '
$ws.Range('E25').Value = '  +4.13%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').Value = '  +5.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.49'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.06'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '161.17'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '34.33'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +8.18%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.13'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0741'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.34%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.83%  '
$ws.Range('E37').Value = '  +2.30%  '
$ws.Range('E38').Value = '  +2.67%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.104'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.37%  '
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.55'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.54%  '
$ws.Range('E41').Value = '  +3.23%  '
$ws.Range('E42').Value = '  +3.98%  '
$ws.Range('D43').Value = '2.056.49'
$ws.Range('E43').Value = '  -0.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.47'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.61%  '
$ws.Range('E45').Value = '  +3.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.13'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.60%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.04'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.84%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.84'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.41%  '
$ws.Range('E49').Value = '  +3.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.20'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.01%  '
$ws.Range('E51').Value = '  +2.80%  '
